$wb = $excel.ActiveWorkbook

# --- Rename Sheet2 to MySheet2 ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Name = "MySheet2"

# --- Sheet3: add a third column ("How I'm feeling about that") ---
$ws3 = $wb.Worksheets.Item("Sheet3")

$ws3.Range("B1").Copy()
$ws3.Range("C1").PasteSpecial(-4122) # xlPasteFormats - copy the header formatting
$ws3.Range("C1").Value = "How I'm feeling about that"
# the longer header text wraps to more lines in the (narrower) new column,
# so the header row grows from 2 lines to 4 lines tall.
$ws3.Rows.Item(1).RowHeight = 57.6

$feelings = @(5,4,2,6,4,3,5,1,2,5,4,2,5,6,1,3,2,1,5,4,6,2,3,4,5,2,3,3,1,4,5)
for ($i = 0; $i -lt $feelings.Length; $i++) {
    $row = $i + 2
    $ws3.Cells.Item($row, 3).Value = $feelings[$i]
}

# --- Sheet3 row 33: move the "^sppooaky" note from A33 to B33 ---
$ws3.Range("A33").ClearContents()
$ws3.Range("B33").Value = "^sppooaky"
$ws3.Rows.Item(33).RowHeight = 28.8

# --- View / selection state ---
# Sheet3 ends up not being the active tab, selection rests on C33.
$ws3.Range("C33").Select()

# MySheet2 becomes the active tab, selection stays on E32.
$ws2.Activate()
$ws2.Range("E32").Select()
